# "Generate Report for Handoff"
#
# The localization-status report moves from "In Translation" to
# "Ready for handoff", and the handoff timestamps are refreshed. The status
# columns/cells (which were holding "In Translation") get widened slightly
# to fit the new, longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: per-language status + last HO xliff generate time ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status
$wsOverview.Range("G2").Value = "2016-08-17 16:38:55" # Latest HO Xliff Generate Date

# --- zh-cn sheet: status + latest handoff datetime ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-17 16:38:50"

# --- de-de sheet: status + latest handoff datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-17 16:38:55"

# --- Widen the status columns so the longer "Ready for handoff" text fits ---
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
